# Added bounded context scope
#
# A new "Bounded Context" column is inserted in front of the existing table
# (old column A "Use Case" shifts to B, ... old column F "Expected Result"
# shifts to G) and every test-case row is tagged with the "Auth" bounded
# context.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A - this shifts the existing columns A:F to B:G
# and keeps every existing value/format intact.
$null = $ws.Columns.Item(1).Insert()

# New column header + the bounded-context value for each data row.
$ws.Range("A1").Value = "Bounded Context"
for ($r = 2; $r -le 6; $r++) {
  $ws.Range("A$r").Value = "Auth"
}

# Give the new header cell (A1) the same look as the rest of the header row.
$null = $ws.Range("B1").Copy()
$null = $ws.Range("A1").PasteSpecial(-4122)

# Column widths: A is brand new, B/C are nudged slightly by the reflow that
# inserting a column triggers. D:G keep their original widths.
$ws.Columns.Item(1).ColumnWidth = 19.2517
$ws.Columns.Item(2).ColumnWidth = 10.4167
$ws.Columns.Item(3).ColumnWidth = 46.0867

# Row heights reflow slightly because of the extra column.
$ws.Rows.Item(1).RowHeight = 27.7
for ($r = 2; $r -le 6; $r++) {
  $ws.Rows.Item($r).RowHeight = 22.85
}

$null = $ws.Range("A8").Select()
